# Excel Basics ( 17-7-2025 )
# Adds a "Shortcut keys / Action" mini-table (cols T:AA, rows 1-7) to Sheet1,
# a new dated row (17-7-2025 / Excel Basics), a small calculation area
# (F8:G11), extends the existing description text + merge, and appends a
# new empty "Sheet2" worksheet after "Sheet1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Row 1 header block: extend the bold/filled header band rightwards
#    and add the "Shortcut keys" / "Action" header cells (cols T..AA).
# ---------------------------------------------------------------------

# N1:S1 used to just be part of the plain bold band; make them match the
# centered bold header look used elsewhere in the row.
$ws.Range("N1:S1").HorizontalAlignment = -4108

# New headers, styled like the existing "Topic"/"Decription" headers
# (bold, filled, centered).
$hdr = $ws.Range("T1:AA1")
$hdr.Font.Bold = $true
$hdr.Interior.Color = $ws.Range("B1").Interior.Color
$hdr.HorizontalAlignment = -4108

$ws.Range("T1").Value = "Shortcut keys"
$ws.Range("T1:U1").Merge()

$ws.Range("V1").Value = "Action"
$ws.Range("V1:W1").Merge()

$ws.Range("X1:Y1").Merge()
$ws.Range("Z1:AA1").Merge()

# ---------------------------------------------------------------------
# 2. Row 2: extend the "Decription" text + merge by one column, and add
#    the first data row of the shortcut-keys table.
# ---------------------------------------------------------------------
$ws.Range("D2:M2").UnMerge()
$ws.Range("D2").Value = "UI of excel, Quick access, Tabs, Ribbons, Namebox , Formual Bar, Rows, Columns, Cell, Workbook and Worksheets"
$ws.Range("D2:N2").Merge()

$ws.Range("T2").Value = "CTRL+Up/Down"
$ws.Range("T2:U2").Merge()

$ws.Range("V2").Value = "Last row and column "
$ws.Range("V2:W2").Merge()

# ---------------------------------------------------------------------
# 3. Row 3 (new): second dated entry + more shortcut-keys data.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "17-7-2025"

$ws.Range("B3").Value = "Excel Basics"
$ws.Range("B3:C3").Merge()

$ws.Range("T3").Value = "Shift F11"
$ws.Range("T3:U3").Merge()

$ws.Range("V3").Value = "add New  Worksheet"
$ws.Range("V3:W3").Merge()

# ---------------------------------------------------------------------
# 4. Row 4 (new): more shortcut-keys data.
# ---------------------------------------------------------------------
$ws.Range("T4").Value = "CTRL pgup/pgdwn"
$ws.Range("T4:U4").Merge()

$ws.Range("V4").Value = "Switch Worksheet"
$ws.Range("V4:W4").Merge()

# ---------------------------------------------------------------------
# 5-7. Rows 5-7 (new): remaining shortcut-keys / workbook actions.
# ---------------------------------------------------------------------
$ws.Range("T5").Value = "CTRL + W"
$ws.Range("V5").Value = "Close Workbook"

$ws.Range("T6").Value = "CTRL+N"
$ws.Range("V6").Value = "New  Workbook"
$ws.Range("V6:W6").Merge()

$ws.Range("T7").Value = "CTRL+O"
$ws.Range("V7").Value = "Open Workbook"
$ws.Range("V7:W7").Merge()

# ---------------------------------------------------------------------
# 8-11. Small calculation scratch area.
# ---------------------------------------------------------------------
$ws.Range("F8").Value = 45
$ws.Range("F9").Formula = "=55+45"
$ws.Range("F11").Formula = "=F8+F9"
$ws.Range("G11").Formula = "=SUM(F8:F9)"

# ---------------------------------------------------------------------
# Selection / view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("M7").Select()

# ---------------------------------------------------------------------
# 12. Add a new blank "Sheet2" after "Sheet1".
# ---------------------------------------------------------------------
$wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("Sheet1")) | Out-Null

# Keep "Sheet1" as the active/visible tab, like in the saved workbook.
$ws.Activate()
$ws.Range("M7").Select()
